# ---------------------------------------------------------------------------
# Co-op Report.docx edit:
#   1) "Information about the Employer" list: fill the blank bullet between
#      "Tim" and "QA" with "Devaraj", then add a new "Gillian" bullet after
#      it that carries the document's "_GoBack" bookmark.
#   2) That bookmark is removed from its old spot inside the Job Description
#      paragraph (a side effect of it being unique/moved, handled below).
#   3) "Conclusions" list: fill the blank trailing bullet with
#      "Automation and QA", then add a new "Corporate environment" bullet
#      after it.
#   4) "Acknowledgments" list: add a new "Devaraj" sub-bullet right after
#      "Venkat".
#
# Edits are applied bottom-to-top (by paragraph index) so that inserting
# paragraphs lower in the document never shifts the index of an
# not-yet-processed target higher up... actually the opposite: targets
# further down the document are handled first so edits made near the top
# of the document (later in this script) don't disturb indices already
# used for targets further down.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 4: Acknowledgments list - add "Devaraj" right after "Venkat"
# ---------------------------------------------------------------------------
$pVenkat = $d.Paragraphs.Item(34)
if ($pVenkat.Range.Text.Trim() -ne "Venkat") {
    throw "Expected paragraph 34 to be 'Venkat', found: " + $pVenkat.Range.Text
}
$pVenkat.Range.InsertParagraphAfter()
$pDevarajAck = $d.Paragraphs.Item(35)
$pDevarajAck.Range.Text = "Devaraj"

# ---------------------------------------------------------------------------
# Change 3: Conclusions list - add "Automation and QA" to the empty trailing
# paragraph, then a new "Corporate environment" paragraph after it.
# ---------------------------------------------------------------------------
$pEmptyConclusions = $d.Paragraphs.Item(31)
if ($pEmptyConclusions.Range.Text.Trim() -ne "") {
    throw "Expected paragraph 31 to be blank, found: " + $pEmptyConclusions.Range.Text
}
$pEmptyConclusions.Range.Text = "Automation and QA"
$pEmptyConclusions.Range.InsertParagraphAfter()
$pCorporate = $d.Paragraphs.Item(32)
$pCorporate.Range.Text = "Corporate environment"

# ---------------------------------------------------------------------------
# Changes 1 & 2: Information about the Employer list - add "Devaraj" to the
# empty paragraph between "Tim" and "QA", then a new "Gillian" paragraph
# after it carrying the "_GoBack" bookmark (moved from its old spot in the
# Job Description paragraph, since a document can only have one bookmark of
# a given name).
# ---------------------------------------------------------------------------
$pEmptyEmployer = $d.Paragraphs.Item(9)
if ($pEmptyEmployer.Range.Text.Trim() -ne "") {
    throw "Expected paragraph 9 to be blank, found: " + $pEmptyEmployer.Range.Text
}
$pEmptyEmployer.Range.Text = "Devaraj"
$pEmptyEmployer.Range.InsertParagraphAfter()
$pGillian = $d.Paragraphs.Item(10)

# Type a temporary trailing character so the bookmark's collapsed insertion
# point is not literally the last character before the paragraph mark
# (inserting a zero-length bookmark exactly there drops it to the start of
# the document); the placeholder is stripped right after the bookmark has
# been anchored.
$pGillian.Range.Text = "GillianX"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$bmPos = $pGillian.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($pGillian.Range.End - 2, $pGillian.Range.End - 1)
$placeholder.Text = ""
